$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue $ws "D2" '60.929.61'
Set-TextValue $ws "E2" '  -2.30%  '
Set-TextValue $ws "D3" '2.420.16'
Set-TextValue $ws "E3" '  -1.49%  '
Set-TextValue $ws "D4" '0.994'
Set-TextValue $ws "E4" '  -0.47%  '
Set-TextValue $ws "D5" '571.12'
Set-TextValue $ws "E5" '  -0.86%  '
Set-TextValue $ws "D6" '140.17'
Set-TextValue $ws "E6" '  -2.93%  '
Set-TextValue $ws "E7" '  +0.12%  '
Set-TextValue $ws "E8" '  -0.78%  '
Set-TextValue $ws "D9" '2.404.68'
Set-TextValue $ws "E9" '  -1.95%  '
Set-TextValue $ws "E10" '  -0.79%  '
Set-TextValue $ws "E11" '  -0.13%  '
Set-TextValue $ws "E12" '  -2.29%  '
Set-TextValue $ws "E13" '  -1.00%  '
Set-TextValue $ws "D14" '26.01'
Set-TextValue $ws "E14" '  -1.28%  '
Set-TextValue $ws "E15" '  -1.62%  '
Set-TextValue $ws "D16" '2.825.65'
Set-TextValue $ws "E16" '  -2.59%  '
Set-TextValue $ws "D17" '60.775.55'
Set-TextValue $ws "E17" '  -2.30%  '
Set-TextValue $ws "D18" '2.405.46'
Set-TextValue $ws "E18" '  -1.80%  '
Set-TextValue $ws "D19" '10.64'
Set-TextValue $ws "D20" '7.42'
Set-TextValue $ws "E20" '  +3.83%  '
Set-TextValue $ws "D21" '323.21'
Set-TextValue $ws "E21" '  -1.60%  '
Set-TextValue $ws "D22" '4.07'
Set-TextValue $ws "E22" '  -1.23%  '
Set-TextValue $ws "E23" '  +0.87%  '
Set-TextValue $ws "E24" '  -0.05%  '
Set-TextValue $ws "D25" '1.89'
Set-TextValue $ws "E25" '  -3.16%  '
Set-TextValue $ws "D26" '64.98'
Set-TextValue $ws "E26" '  -1.09%  '
Set-TextValue $ws "D27" '8.53'
Set-TextValue $ws "E27" '  -7.00%  '
Set-TextValue $ws "D28" '577.41'
Set-TextValue $ws "E28" '  -2.81%  '
Set-TextValue $ws "D29" '2.515.81'
Set-TextValue $ws "D30" '0.0₃0923'
Set-TextValue $ws "E30" '  -3.62%  '
Set-TextValue $ws "E31" '  -1.01%  '
Set-TextValue $ws "E32" '  -5.46%  '
Set-TextValue $ws "D33" '1.85'
Set-TextValue $ws "E33" '  -1.70%  '
Set-TextValue $ws "E34" '  -2.88%  '
Set-TextValue $ws "E35" '  +0.11%  '
Set-TextValue $ws "D36" '4.67'
Set-TextValue $ws "E36" '  -4.87%  '
Set-TextValue $ws "E37" '  -2.95%  '
Set-TextValue $ws "E38" '  -2.42%  '
Set-TextValue $ws "D39" '149.49'
Set-TextValue $ws "E39" '  -1.66%  '
Set-TextValue $ws "D40" '18.29'
Set-TextValue $ws "E40" '  -0.89%  '
Set-TextValue $ws "D41" '5.17'
Set-TextValue $ws "E41" '  -3.75%  '
Set-TextValue $ws "E42" '  +0.06%  '
Set-TextValue $ws "E43" '  -3.21%  '
Set-TextValue $ws "D44" '41.12'
Set-TextValue $ws "E44" '  -3.35%  '
Set-TextValue $ws "D45" '2.35'
Set-TextValue $ws "E45" '  -5.19%  '
Set-TextValue $ws "D46" '0.0₆0277'
Set-TextValue $ws "E46" '  +12.86%  '
Set-TextValue $ws "D47" '141.38'
Set-TextValue $ws "E47" '  -0.76%  '
Set-TextValue $ws "E48" '  -3.07%  '
Set-TextValue $ws "E49" '  -2.64%  '
Set-TextValue $ws "D50" '19.63'
Set-TextValue $ws "E50" '  -0.71%  '
Set-TextValue $ws "D51" '0.0507'
Set-TextValue $ws "E51" '  -3.05%  '
